$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '89.459.41'
$ws.Range("E2").Value = '  -1.70%  '
$ws.Range("D3").Value = '3.075.49'
$ws.Range("E3").Value = '  -2.63%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '616.26'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.06'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.94%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.362'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").Value = '3.074.07'
$ws.Range("E10").Value = '  -2.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.710'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.198'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000250'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.54%  '
$ws.Range("D15").Value = '89.364.83'
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.35'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.71%  '
$ws.Range("D17").Value = '3.641.50'
$ws.Range("E17").Value = '  -2.60%  '
$ws.Range("D18").Value = '3.084.21'
$ws.Range("E18").Value = '  -3.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000211'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '430.44'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -9.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '86.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.79%  '
$ws.Range("D28").Value = '3.250.94'
$ws.Range("E28").Value = '  -2.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +14.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.99'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.156'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.196'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -15.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.52'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.149'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.08'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '492.59'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.62'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.88'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.25'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.64'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +54.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0892'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.396'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '152.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.84'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.671'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.30'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.30%  '
